# Support initializing multiple variables in one line
# Updates underlying benchmark numbers on the "Concise" sheet; dependent
# formulas (columns J:N) and the chart caches that reference this sheet
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concise")

# Row 3 (Python baseline row)
$ws.Range("D3").Value = 5600
$ws.Range("F3").Value = 2687
$ws.Range("G3").Value = 682

# Row 4 (Bau)
$ws.Range("C4").Value = 228
$ws.Range("D4").Value = 5513
$ws.Range("E4").Value = 2064
$ws.Range("F4").Value = 2805
$ws.Range("G4").Value = 641

# Row 5 (Swift)
$ws.Range("D5").Value = 6564
$ws.Range("E5").Value = 2540
$ws.Range("F5").Value = 3195
$ws.Range("G5").Value = 820

# Row 7 (C)
$ws.Range("C7").Value = 306
$ws.Range("D7").Value = 7273
$ws.Range("E7").Value = 2606
$ws.Range("F7").Value = 3473
$ws.Range("G7").Value = 1126

# Row 8 (Go)
$ws.Range("C8").Value = 334
$ws.Range("D8").Value = 7063
$ws.Range("E8").Value = 2775
$ws.Range("F8").Value = 3213
$ws.Range("G8").Value = 1016

# Row 9 (Java)
$ws.Range("D9").Value = 7928
$ws.Range("F9").Value = 3601
$ws.Range("G9").Value = 1004

$excel.CalculateFullRebuild()

# Restore the sheet's active selection/scroll position
$ws.Activate()
$ws.Range("L11").Select()
